$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Z2").Value = 4.33
# Row 3
$ws.Range("G3").Value = 2.25
$ws.Range("H3").Value = 3.1
$ws.Range("P3").Value = 1.53
$ws.Range("Q3").Value = 2.38
$ws.Range("AC3").Value = 67
$ws.Range("AD3").Value = 451
$ws.Range("AF3").Value = 15
# Row 4
$ws.Range("J4").Value = 1.1
$ws.Range("K4").Value = 7
# Row 8
$ws.Range("L8").Value = 1.5
$ws.Range("M8").Value = 2.5
# Row 10
$ws.Range("G10").Value = 1.9
$ws.Range("H10").Value = 3.2
$ws.Range("J10").Value = 1.1
$ws.Range("K10").Value = 7
$ws.Range("U10").Value = 7.5
$ws.Range("AA10").Value = 6.5
$ws.Range("AB10").Value = 21
$ws.Range("AC10").Value = 81
$ws.Range("AG10").Value = 15
# Row 11
$ws.Range("G11").Value = 1.73
$ws.Range("H11").Value = 3.3
$ws.Range("I11").Value = 4.75
$ws.Range("K11").Value = 7.5
$ws.Range("N11").Value = 2.25
$ws.Range("O11").Value = 1.62
$ws.Range("R11").Value = 2.1
$ws.Range("S11").Value = 1.67
$ws.Range("T11").Value = 5.5
$ws.Range("U11").Value = 7.5
$ws.Range("W11").Value = 13
$ws.Range("AE11").Value = 11
$ws.Range("AF11").Value = 23
$ws.Range("AG11").Value = 17
$ws.Range("AJ11").Value = 51
# Row 12
$ws.Range("G12").Value = 1.44
$ws.Range("H12").Value = 4.1
$ws.Range("I12").Value = 6.5
$ws.Range("N12").Value = 1.8
$ws.Range("O12").Value = 2
$ws.Range("R12").Value = 1.91
$ws.Range("S12").Value = 1.8
$ws.Range("X12").Value = 12
$ws.Range("Y12").Value = 26
$ws.Range("AE12").Value = 17
$ws.Range("AG12").Value = 21
$ws.Range("AH12").Value = 81
# Row 13
$ws.Range("P13").Value = 1.5
# Row 14
$ws.Range("J14").Value = 1.05
$ws.Range("K14").Value = 11
$ws.Range("L14").Value = 1.29
$ws.Range("M14").Value = 3.5
$ws.Range("N14").Value = 1.98
$ws.Range("O14").Value = 1.83
$ws.Range("P14").Value = 1.44
$ws.Range("Q14").Value = 2.63
# Row 15
$ws.Range("N15").Value = 1.85
$ws.Range("O15").Value = 2
$ws.Range("P15").Value = 1.36
# Row 18
$ws.Range("L18").Value = 1.3
$ws.Range("M18").Value = 3.4
# Row 19
$ws.Range("J19").Value = 1.04
$ws.Range("K19").Value = 13
$ws.Range("L19").Value = 1.29
$ws.Range("M19").Value = 3.5
$ws.Range("N19").Value = 1.9
$ws.Range("O19").Value = 1.95
# Row 20
$ws.Range("L20").Value = 1.29
$ws.Range("M20").Value = 3.5
$ws.Range("N20").Value = 1.95
$ws.Range("O20").Value = 1.85
# Row 21
$ws.Range("G21").Value = 3.4
$ws.Range("I21").Value = 2.3
$ws.Range("P21").Value = 1.57
$ws.Range("Q21").Value = 2.25
$ws.Range("R21").Value = 2.1
$ws.Range("S21").Value = 1.67
$ws.Range("W21").Value = 41
$ws.Range("X21").Value = 34
$ws.Range("AD21").Value = 1000
$ws.Range("AE21").Value = 6
$ws.Range("AH21").Value = 21
$ws.Range("AI21").Value = 21

Write-Host "Applied odds updates to Sheet1"
